# HSA VISTA GENERAL - Actualización Anyi 8 de mayo
# The "ANYI AGUIRRE" sheet had its second data row (row 2) removed -
# that record (CNE-E-DG-2025-004014 / QUEJA CONTRA EL SEÑOR DANIEL BECERRA ...)
# was deleted, shifting the remaining rows up and making this the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ANYI AGUIRRE")

# Make it the active sheet (matches activeTab pointing at ANYI AGUIRRE)
$ws.Select()

# Select and delete the whole second row, shifting the rows below it up
$row = $ws.Rows.Item(2)
$row.Select()
$row.Delete()

# After the delete, Excel leaves the (new) row 2 selected
$ws.Rows.Item(2).Select()
